# GSP362_TeamA_BugTracking.xlsx - "Update Bugs and Build Settings"
#  - Updated the bug sheet.
#  - Added the other scenes to the build settings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 6: "No defined end game" -> mark completed, add Solution note ---
$ws.Range("C6").Value = "completed"

# --- Row 7: "Solo play does not work" -> mark completed, add Solution note ---
# (write E7 before E6 so new shared strings land in the same order Excel produced them)
$ws.Range("E7").Value = "Added collision the the map prefab and functionality to make the game end upon this collision."
$ws.Range("C7").Value = "completed"
$ws.Range("E6").Value = "Added support for passing number of players from the menu to the game scene."

# --- Row 9: new bug entry ---
$ws.Range("A9").Value = "Resources are editor only."
$ws.Range("B9").Value = "high"
$ws.Range("C9").Value = "in progress"
$ws.Range("D9").Value = "Resources are loaded with editor functions. This isn't allowed in a built game."

# --- Row heights grow because the new notes wrap onto multiple lines ---
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 30

# --- Selection moved to D9 ---
$ws.Activate()
$ws.Range("D9").Select() | Out-Null
